# Add a "schema" field to the BenchmarkResult schema-definition sheets
# (and extend the "id" validation drop-down list to include every schema id).
#
# Each of the two sheets currently ends its header row with column I
# holding "summary". We insert a new column at I so the old I ("summary")
# shifts right to J, and the freshly-opened column I gets the label
# "schema" - matching the xl diff (dimension A1:I1 -> A1:J1, new I1 cell
# "schema", old I1 "summary" becomes J1).

$wb = $excel.ActiveWorkbook

$targetSheets = @("BenchmarkResult", "BenchmarkResult1")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new blank column before column I (9) - this shifts the
    # existing I1 ("summary") to J1 and widens the sheet dimension from
    # A1:I1 to A1:J1 automatically.
    $ws.Columns.Item(9).Insert()

    # Label the newly inserted column.
    $ws.Range("I1").Value = "schema"
}

# The "BenchmarkResult" sheet's `id` column (B) is a drop-down list of
# every valid schema id. Extend it with the newly added schema ids.
$newIdList = '"1a,1b,1c,1d,2a,2b,2c,2d,3a,4a,4b,4c,4d,4e,4f,4g,4h,5a,5b,6a,6b,7a,7b,7c,8a,8b,8c"'
$wsBenchmarkResult = $wb.Worksheets.Item("BenchmarkResult")
$wsBenchmarkResult.Range("B2:B1048576").Validation.Modify(3, 1, 1, $newIdList)
